# Weekly update: a new daily record for "Arveja Verde" (Hortaliza) was
# reported for Vega Central Mapocho de Santiago, and the existing history
# (previously rows 60-102) shifts down by one row to make room for it.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row above the current row 60, pushing the existing
# rows 60-102 down to 61-103 (dimension grows from R102 to R103).
$ws.Rows("60").Insert(-4121)

# Populate the newly inserted row 60 with the new record.
$ws.Range("A60").Value = 9
$ws.Range("B60").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C60").Value = "Metropolitana"
$ws.Range("D60").Value = 44596
$ws.Range("E60").Value = 13
$ws.Range("F60").Value = 100112022
$ws.Range("G60").Value = "Arveja Verde"
$ws.Range("H60").Value = "Sin especificar"
$ws.Range("I60").Value = "Primera"
$ws.Range("J60").Value = 30
$ws.Range("K60").Value = 28000
$ws.Range("L60").Value = 28000
$ws.Range("M60").Value = 28000
$ws.Range("N60").Value = "$/saco 25 kilos"
$ws.Range("O60").Value = "Carahue"
$ws.Range("P60").Value = 1120
$ws.Range("Q60").Value = 25
$ws.Range("R60").Value = "Hortaliza"
